$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ALNS staff schedule output (A2:AC10)
$data = @(
    @("Staff_1", "A1", "A1", "A1", "DO", "A1", "M1", "M3", "A1", "M3", "DO", "A1", "A1", "M1", "M1", "DO", "M3", "A1", "A1", "A1", "M1", "A1", "M3", "M1", "M1", "A1", "A1", "M1", "DO"),
    @("Staff_2", "A1", "M2", "M2", "A2", "M2", "M3", "DO", "A1", "DO", "M2", "M2", "M2", "M3", "M3", "A1", "A1", "DO", "M2", "M1", "M3", "M3", "M2", "A1", "DO", "M2", "M3", "M3", "M2"),
    @("Staff_3", "M1", "DO", "A1", "M1", "M3", "M3", "M1", "M1", "M1", "A1", "M2", "M1", "M3", "DO", "M1", "M1", "M2", "DO", "M1", "A1", "M3", "M1", "M1", "DO", "M1", "M3", "A1", "M1"),
    @("Staff_4", "M2", "M1", "DO", "M2", "M3", "M3", "A2", "M1", "M2", "M2", "DO", "M1", "M3", "A2", "A1", "M2", "M1", "DO", "M3", "M2", "M1", "M3", "DO", "M1", "M2", "M2", "M2", "A1"),
    @("Staff_5", "M1", "M1", "M1", "DO", "M2", "M3", "M1", "M3", "A1", "M1", "M1", "DO", "M1", "M1", "A2", "A2", "DO", "M3", "M1", "M2", "M1", "A2", "DO", "A2", "M1", "M1", "M3", "A2"),
    @("Staff_6", "DO", "M3", "A2", "A2", "A2", "A1", "A2", "DO", "A2", "M3", "A1", "A2", "A2", "A2", "M3", "DO", "A1", "A1", "A1", "A1", "A1", "A1", "A1", "A1", "A1", "DO", "A1", "M3"),
    @("Staff_7", "A2", "M3", "A2", "A2", "A1", "M2", "DO", "DO", "A2", "A1", "A2", "A1", "M3", "M2", "A2", "M3", "A1", "A1", "DO", "M1", "M1", "M1", "A2", "A2", "A2", "A1", "M3", "DO"),
    @("Staff_8", "M3", "M1", "M2", "M2", "DO", "A2", "A1", "M3", "DO", "M2", "M1", "M2", "A1", "A2", "M1", "DO", "M2", "M1", "M3", "A2", "A1", "DO", "M2", "M2", "M2", "M3", "A1", "A1"),
    @("Staff_9", "DO", "A1", "M2", "M2", "M3", "A2", "M2", "A2", "M2", "A1", "DO", "A1", "A2", "M3", "M1", "M3", "M1", "M2", "A1", "DO", "A2", "DO", "M1", "M2", "M2", "A2", "M3", "M3")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

